$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "65.874.08"; E = "  +1.09%  " },
    @{ Row = 3;  D = "3.198.74";  E = $null },
    @{ Row = 4;  D = $null;       E = "  -0.05%  " },
    @{ Row = 5;  D = "599.55";    E = "  +3.60%  " },
    @{ Row = 6;  D = "152.95";    E = "  +1.31%  " },
    @{ Row = 7;  D = $null;       E = "  +0.02%  " },
    @{ Row = 8;  D = "3.195.07";  E = "  +0.81%  " },
    @{ Row = 9;  D = $null;       E = "  +0.29%  " },
    @{ Row = 10; D = "0.159";     E = "  -1.48%  " },
    @{ Row = 11; D = "6.07";      E = "  -2.23%  " },
    @{ Row = 12; D = $null;       E = "  +1.51%  " },
    @{ Row = 13; D = $null;       E = "  -0.49%  " },
    @{ Row = 14; D = "39.47";     E = "  +5.23%  " },
    @{ Row = 15; D = "3.723.45";  E = "  +1.09%  " },
    @{ Row = 16; D = "7.47";      E = "  +4.03%  " },
    @{ Row = 17; D = "65.993.76"; E = "  +1.20%  " },
    @{ Row = 18; D = "3.201.77";  E = "  +1.13%  " },
    @{ Row = 19; D = $null;       E = "  -0.03%  " },
    @{ Row = 20; D = "510.22";    E = "  -0.28%  " },
    @{ Row = 21; D = "15.40";     E = "  +3.61%  " },
    @{ Row = 22; D = $null;       E = "  +1.68%  " },
    @{ Row = 23; D = "8.17";      E = "  +4.40%  " },
    @{ Row = 24; D = "15.34";     E = "  +0.22%  " },
    @{ Row = 25; D = "84.79";     E = "  -0.33%  " },
    @{ Row = 26; D = "0.998";     E = "  +0.01%  " },
    @{ Row = 27; D = "9.28";      E = "  +2.22%  " },
    @{ Row = 28; D = $null;       E = "  +2.65%  " },
    @{ Row = 29; D = $null;       E = "  +3.54%  " },
    @{ Row = 30; D = $null;       E = "  +8.59%  " },
    @{ Row = 31; D = "2.87";      E = "  +1.48%  " },
    @{ Row = 32; D = "28.06";     E = $null },
    @{ Row = 33; D = $null;       E = "  +1.84%  " },
    @{ Row = 34; D = $null;       E = "  +0.01%  " },
    @{ Row = 35; D = $null;       E = "  -0.79%  " },
    @{ Row = 36; D = "54.99";     E = "  -1.49%  " },
    @{ Row = 37; D = "0.0906";    E = "  +0.44%  " },
    @{ Row = 38; D = "484.46";    E = "  +2.08%  " },
    @{ Row = 39; D = "0.0419";    E = "  -0.16%  " },
    @{ Row = 40; D = $null;       E = "  -3.92%  " },
    @{ Row = 41; D = "8.88";      E = "  +2.36%  " },
    @{ Row = 42; D = "0.302";     E = "  +5.60%  " },
    @{ Row = 43; D = $null;       E = "  +1.87%  " },
    @{ Row = 44; D = "2.949.65";  E = "  -4.01%  " },
    @{ Row = 45; D = $null;       E = "  +6.77%  " },
    @{ Row = 46; D = $null;       E = "  -0.95%  " },
    @{ Row = 47; D = "28.47";     E = "  -2.54%  " },
    @{ Row = 49; D = $null;       E = "  +0.87%  " },
    @{ Row = 50; D = $null;       E = "  +2.15%  " },
    @{ Row = 51; D = "2.59";      E = "  +4.16%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Range("D$($u.Row)")
        # Force text so values like "599.55" or "0.159" are not
        # reinterpreted as numbers, then restore the default style
        # so no extra formatting is left behind on the cell.
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $ws.Range("E$($u.Row)").Value = $u.E
    }
}
